$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "69.404.64"
$ws.Range("E2").Value = "  +1.66%  "

# Row 3
$ws.Range("D3").Value = "3.902.30"
$ws.Range("E3").Value = "  +0.38%  "

# Row 4
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").Value = "'531.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +9.90%  "

# Row 6
$ws.Range("D6").Value = "'144.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.62%  "

# Row 7
$ws.Range("D7").Value = "'0.613"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.25%  "

# Row 8
$ws.Range("D8").Value = "'0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.09%  "

# Row 9
$ws.Range("D9").Value = "'0.721"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.50%  "

# Row 10
$ws.Range("E10").Value = "  -1.53%  "

# Row 11
$ws.Range("D11").Value = "'0.0000335"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -4.72%  "

# Row 12
$ws.Range("D12").Value = "'42.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.74%  "

# Row 13
$ws.Range("D13").Value = "4.524.56"
$ws.Range("E13").Value = "  +0.39%  "

# Row 14
$ws.Range("D14").Value = "'10.29"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.58%  "

# Row 15
$ws.Range("D15").Value = "3.909.95"
$ws.Range("E15").Value = "  -0.08%  "

# Row 16
$ws.Range("D16").Value = "'14.03"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.06%  "

# Row 17
$ws.Range("E17").Value = "  -1.33%  "

# Row 18
$ws.Range("E18").Value = "  +6.93%  "

# Row 19
$ws.Range("D19").Value = "'19.83"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.57%  "

# Row 20
$ws.Range("D20").Value = "69.435.69"
$ws.Range("E20").Value = "  +1.71%  "

# Row 21
$ws.Range("D21").Value = "'425.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.98%  "

# Row 22
$ws.Range("E22").Value = "  -3.96%  "

# Row 23
$ws.Range("E23").Value = "  -3.87%  "

# Row 24
$ws.Range("D24").Value = "'88.12"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.80%  "

# Row 25
$ws.Range("D25").Value = "'4.04"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +9.33%  "

# Row 26
$ws.Range("D26").Value = "'11.42"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -9.06%  "

# Row 27
$ws.Range("D27").Value = "'10.59"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.36%  "

# Row 28
$ws.Range("D28").Value = "'36.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.93%  "

# Row 29
$ws.Range("D29").Value = "'690.61"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -4.03%  "

# Row 30
$ws.Range("D30").Value = "'13.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.05%  "

# Row 31
$ws.Range("E31").Value = "  -2.73%  "

# Row 32
$ws.Range("E32").Value = "  -2.74%  "

# Row 33
$ws.Range("D33").Value = "'68.69"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.07%  "

# Row 34
$ws.Range("D34").Value = "0.0₃0872"
$ws.Range("E34").Value = "  -0.18%  "

# Row 35
$ws.Range("D35").Value = "'0.433"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +8.39%  "

# Row 36
$ws.Range("E36").Value = "  -1.61%  "

# Row 37
$ws.Range("D37").Value = "'40.02"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.04%  "

# Row 38
$ws.Range("D38").Value = "'0.149"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +2.46%  "

# Row 39
$ws.Range("D39").Value = "'1.00"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.09%  "

# Row 40
$ws.Range("B40").Value = "ThetaToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D40").Value = "'3.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +8.47%  "

# Row 41
$ws.Range("B41").Value = "FirstDigitalUSD"
$ws.Range("C41").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.16%  "

# Row 42
$ws.Range("D42").Value = "'0.0484"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.24%  "

# Row 43
$ws.Range("E43").Value = "  +7.88%  "

# Row 44
$ws.Range("D44").Value = "'2.78"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -7.07%  "

# Row 45
$ws.Range("D45").Value = "'3.40"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.84%  "

# Row 46
$ws.Range("D46").Value = "'0.140"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.05%  "

# Row 47
$ws.Range("D47").Value = "'0.000285"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +15.50%  "

# Row 48
$ws.Range("E48").Value = "  +7.20%  "

# Row 49
$ws.Range("D49").Value = "0.0₆0344"
$ws.Range("E49").Value = "  -4.60%  "

# Row 50
$ws.Range("B50").Value = "Maker"
$ws.Range("C50").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D50").Value = "2.747.63"
$ws.Range("E50").Value = "  +15.27%  "

# Row 51
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").Value = "'145.80"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.06%  "
